$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.114.22'
$ws.Range('E2').Value = '  -0.27%  '
$ws.Range('D3').Value = '1.875.98'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''313.51'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '''0.5057'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '''0.3836'
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('D9').Value = '''0.08562'
$ws.Range('E9').Value = '  -8.33%  '
$ws.Range('E10').Value = '  -2.37%  '
$ws.Range('D11').Value = '''41.35'
$ws.Range('E11').Value = '  -1.53%  '
$ws.Range('D12').Value = '''6.306'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').Value = '''20.66'
$ws.Range('E13').Value = '  -1.48%  '
$ws.Range('D14').Value = '1.876.16'
$ws.Range('E14').Value = '  -2.13%  '
$ws.Range('D15').Value = '''7.207'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').Value = '''1.003'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('E17').Value = '  -2.56%  '
$ws.Range('D18').Value = '''91.03'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D19').Value = '''0.06629'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').Value = '''18.10'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('D22').Value = '''6.095'
$ws.Range('E22').Value = '  -2.18%  '
$ws.Range('D23').Value = '28.135.22'
$ws.Range('D24').Value = '''11.39'
$ws.Range('E24').Value = '  -1.44%  '
$ws.Range('D25').Value = '''2.260'
$ws.Range('E25').Value = '  -2.90%  '
$ws.Range('D26').Value = '''2.591'
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = '2.090.70'
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').Value = '  -2.18%  '
$ws.Range('D29').Value = '''156.97'
$ws.Range('E29').Value = '  -0.79%  '
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('D31').Value = '''0.1058'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('E32').Value = '  -4.04%  '
$ws.Range('D33').Value = '''5.618'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').Value = '''3.592'
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('D35').Value = '''9.616'
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('D36').Value = '''0.02453'
$ws.Range('E36').Value = '  +0.62%  '
$ws.Range('D37').Value = '''0.06589'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('E38').Value = '  -1.56%  '
$ws.Range('D39').Value = '''1.214'
$ws.Range('D40').Value = '''1.245'
$ws.Range('E40').Value = '  -2.87%  '
$ws.Range('D41').Value = '''0.6378'
$ws.Range('E41').Value = '  -2.39%  '
$ws.Range('D42').Value = '''11.45'
$ws.Range('E42').Value = '  -0.97%  '
$ws.Range('D43').Value = '''4.891'
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''13.20'
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.6005'
$ws.Range('E45').Value = '  -1.91%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = '''1.285'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').Value = '''3.675'
$ws.Range('E47').Value = '  -1.29%  '
$ws.Range('B48').Value = 'EOS'
$ws.Range('C48').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D48').Value = '''1.225'
$ws.Range('E48').Value = '  +2.94%  '
$ws.Range('E49').Value = '  -2.07%  '
$ws.Range('B50').Value = 'Quant'
$ws.Range('C50').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D50').Value = '''121.61'
$ws.Range('E50').Value = '  -0.69%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''80.66'
$ws.Range('E51').Value = '  +2.86%  '
